$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update totals in row 2 ("мелочный товар") ---
$ws.Range("AL2").Value = 5
$ws.Range("AP2").Value = 271

# --- Update totals in row 9 ("мелкий товар") ---
$ws.Range("N9").Value = 10
$ws.Range("AP9").Value = 31

# --- Re-run of text analysis shuffled counts for rows 28-41 ---
# Row 28
$ws.Range("A28").Value = "внутренний товар"
$ws.Range("I28").Value = ""
$ws.Range("AB28").Value = 1
$ws.Range("AH28").Value = 1
$ws.Range("AM28").Value = 1
$ws.Range("AN28").Value = ""
# Row 29
$ws.Range("A29").Value = "суровский товар"
$ws.Range("N29").Value = ""
$ws.Range("AB29").Value = ""
$ws.Range("AH29").Value = ""
$ws.Range("AM29").Value = ""
$ws.Range("AO29").Value = 4
# Row 31
$ws.Range("A31").Value = "питейный припасы"
$ws.Range("I31").Value = 1
$ws.Range("L31").Value = ""
$ws.Range("N31").Value = 1
$ws.Range("AI31").Value = ""
# Row 32
$ws.Range("A32").Value = "медный товар"
$ws.Range("L32").Value = 1
$ws.Range("R32").Value = ""
$ws.Range("AE32").Value = ""
$ws.Range("AI32").Value = 1
$ws.Range("AN32").Value = 2
$ws.Range("AP32").Value = 4
# Row 33
$ws.Range("A33").Value = "привозный товар"
$ws.Range("O33").Value = ""
$ws.Range("R33").Value = 1
$ws.Range("AE33").Value = 1
$ws.Range("AI33").Value = ""
# Row 34
$ws.Range("A34").Value = "оловянный товар"
$ws.Range("O34").Value = 1
$ws.Range("AI34").Value = 1
$ws.Range("AN34").Value = 1
$ws.Range("AO34").Value = ""
# Row 39
$ws.Range("A39").Value = "меховой товар"
$ws.Range("R39").Value = 1
$ws.Range("AK39").Value = ""
# Row 40
$ws.Range("A40").Value = "рукодельный товар"
$ws.Range("R40").Value = ""
$ws.Range("AN40").Value = 1
# Row 41
$ws.Range("A41").Value = "надлежащий товар"
$ws.Range("AK41").Value = 1
$ws.Range("AN41").Value = ""
